$win = $excel.ActiveWindow
$win.Zoom = 115
Write-Host "done"
